$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 21.166666
$ws.Range("I8").Value = 21.166666
$ws.Range("K8").Value = 63.499998
$ws.Range("M8").Value = 75.50000199999999
$ws.Range("H28").Value = 1195.3
$ws.Range("I28").Value = 1287.25
$ws.Range("J28").Value = 827.5
$ws.Range("K28").Value = 1287.25
$ws.Range("L28").Value = 827.5
$ws.Range("M28").Value = -802.25
$ws.Range("N28").Value = -1797.5
$ws.Range("H41").Value = 643
$ws.Range("I41").Value = 643
$ws.Range("K41").Value = 643
$ws.Range("M41").Value = -203
$ws.Range("H70").Value = 3449.9167
$ws.Range("J70").Value = 3379.8
$ws.Range("L70").Value = 10139.4
$ws.Range("N70").Value = -10679.4
$ws.Range("H73").Value = 3449.9167
$ws.Range("J73").Value = 3379.8
$ws.Range("L73").Value = 10139.4
$ws.Range("N73").Value = -12011.4
$ws.Range("H88").Value = 2005
$ws.Range("I88").Value = 1685.25
$ws.Range("J88").Value = 2324.75
$ws.Range("K88").Value = 1685.25
$ws.Range("L88").Value = 2324.75
$ws.Range("M88").Value = -1279.25
$ws.Range("N88").Value = -3136.75
$ws.Range("H91").Value = 2005
$ws.Range("I91").Value = 1685.25
$ws.Range("J91").Value = 2324.75
$ws.Range("K91").Value = 1685.25
$ws.Range("L91").Value = 2324.75
$ws.Range("M91").Value = -281.25
$ws.Range("N91").Value = -5132.75
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 16504
$ws.Range("I38").Value = 4840
$ws.Range("K38").Value = 4840
$ws.Range("M38").Value = -4373
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H96").Value = 35377.285
$ws.Range("J96").Value = 35377.285
$ws.Range("L96").Value = 35377.285
$ws.Range("N96").Value = -40869.285
$ws.Range("H97").Value = 1219.8636
$ws.Range("I97").Value = 935.3889
$ws.Range("K97").Value = 935.3889
$ws.Range("M97").Value = -439.3889
$ws.Range("H110").Value = 742.75
$ws.Range("I110").Value = 715.6
$ws.Range("J110").Value = 762.1429000000001
$ws.Range("K110").Value = 715.6
$ws.Range("L110").Value = 762.1429000000001
$ws.Range("M110").Value = 1329.4
$ws.Range("N110").Value = -4852.1429

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1900
$ws.Range("I20").Value = 1900
$ws.Range("K20").Value = 1900
$ws.Range("M20").Value = -1653
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H86").Value = 937.6667
$ws.Range("I86").Value = 925.2
$ws.Range("K86").Value = 925.2
$ws.Range("M86").Value = 197.8
$ws.Range("H89").Value = 937.6667
$ws.Range("I89").Value = 925.2
$ws.Range("K89").Value = 4626
$ws.Range("M89").Value = 990
$ws.Range("H94").Value = 2544.1428
$ws.Range("I94").Value = 1984.8334
$ws.Range("J94").Value = 2963.625
$ws.Range("K94").Value = 1984.8334
$ws.Range("L94").Value = 2963.625
$ws.Range("M94").Value = -1533.8334
$ws.Range("N94").Value = -3865.625

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 348.85715
$ws.Range("I2").Value = 60.5
$ws.Range("J2").Value = 464.2
$ws.Range("K2").Value = 363
$ws.Range("L2").Value = 2785.2
$ws.Range("M2").Value = -250
$ws.Range("N2").Value = -3011.2
$ws.Range("H11").Value = 282.75
$ws.Range("I11").Value = 282.75
$ws.Range("K11").Value = 848.25
$ws.Range("M11").Value = -708.25
$ws.Range("H63").Value = 4093.2
$ws.Range("J63").Value = 5225
$ws.Range("L63").Value = 15675
$ws.Range("N63").Value = -17173
$ws.Range("H66").Value = 4093.2
$ws.Range("J66").Value = 5225
$ws.Range("L66").Value = 47025
$ws.Range("N66").Value = -54513
$ws.Range("H95").Value = 15000
$ws.Range("J95").Value = 15000
$ws.Range("L95").Value = 45000
$ws.Range("N95").Value = -49118
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H113").Value = 280
$ws.Range("I113").Value = 280
$ws.Range("K113").Value = 840
$ws.Range("M113").Value = 1330
$ws.Range("H121").Value = 363.5
$ws.Range("I121").Value = 363.5
$ws.Range("K121").Value = 1090.5
$ws.Range("M121").Value = 219.5
$ws.Range("H137").Value = 6666.6665
$ws.Range("J137").Value = 6666.6665
$ws.Range("L137").Value = 19999.9995
$ws.Range("N137").Value = -30199.9995

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 280408.75
$ws.Range("I11").Value = 350235.94
$ws.Range("J11").Value = 1100
$ws.Range("K11").Value = 350235.94
$ws.Range("L11").Value = 1100
$ws.Range("M11").Value = -350096.94
$ws.Range("N11").Value = -1378
$ws.Range("H53").Value = 46000
$ws.Range("J53").Value = 46000
$ws.Range("L53").Value = 46000
$ws.Range("N53").Value = -47262
$ws.Range("H55").Value = 38500
$ws.Range("J55").Value = 38500
$ws.Range("L55").Value = 38500
$ws.Range("N55").Value = -39154
$ws.Range("H80").Value = 2984.3333
$ws.Range("J80").Value = 3476.75
$ws.Range("L80").Value = 3476.75
$ws.Range("N80").Value = -5472.75
$ws.Range("H83").Value = 2984.3333
$ws.Range("J83").Value = 3476.75
$ws.Range("L83").Value = 17383.75
$ws.Range("N83").Value = -27367.75
$ws.Range("H97").Value = 1199.4445
$ws.Range("J97").Value = 2002.5
$ws.Range("L97").Value = 2002.5
$ws.Range("N97").Value = -2994.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1388.9286
$ws.Range("I82").Value = 1137.625
$ws.Range("J82").Value = 1724
$ws.Range("K82").Value = 1137.625
$ws.Range("L82").Value = 1724
$ws.Range("M82").Value = -776.625
$ws.Range("N82").Value = -2446
$ws.Range("H85").Value = 1388.9286
$ws.Range("I85").Value = 1137.625
$ws.Range("J85").Value = 1724
$ws.Range("K85").Value = 1137.625
$ws.Range("L85").Value = 1724
$ws.Range("M85").Value = 110.375
$ws.Range("N85").Value = -4220

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 953.75
$ws.Range("I13").Value = 953.75
$ws.Range("K13").Value = 953.75
$ws.Range("M13").Value = -813.75
$ws.Range("H81").Value = 2289.5833
$ws.Range("I81").Value = 2307.7273
$ws.Range("J81").Value = 2090
$ws.Range("K81").Value = 4615.4546
$ws.Range("L81").Value = 4180
$ws.Range("M81").Value = -3554.4546
$ws.Range("N81").Value = -6302
$ws.Range("H84").Value = 2289.5833
$ws.Range("I84").Value = 2307.7273
$ws.Range("J84").Value = 2090
$ws.Range("K84").Value = 23077.273
$ws.Range("L84").Value = 20900
$ws.Range("M84").Value = -17773.273
$ws.Range("N84").Value = -31508
